{"js": "// Make the \"Prevention methods:\" heading bold.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet heading = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Prevention methods:\") {\n    heading = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!heading) {\n  throw new Error('Paragraph \"Prevention methods:\" not found.');\n}\n\nheading.font.bold = true;\nawait context.sync();\n", "ps1": "# Make the \"Prevention methods:\" heading bold.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7, [char]10)\n    if ($text -eq \"Prevention methods:\") {\n        $p.Range.Bold = 1\n        break\n    }\n}\n"}
